$d = $word.ActiveDocument

# The text "3. Usuário do Sistema clica no botão 'Calcular Desconto!' bs 18"
# appears as the 3rd step of several Alternative Flow sections (AF[4],
# AF[5], AF[6], AF[7]). The commit ("Melhoria do fluxo alternativo 7")
# removes only the copy that belongs to AF[7] — the last occurrence in
# the document, directly after the MSG002 line. Find every occurrence,
# remember the last one's range, then delete that whole paragraph
# (its run, the empty lead-in run, and the bookmarkEnd it carries).
$searchText = "3. Usuário do Sistema clica no botão 'Calcular Desconto!' bs 18"

$rng = $d.Content
$lastStart = -1
$lastEnd = -1
while ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastStart = $rng.Start
    $lastEnd = $rng.End
    $rng.Collapse(0)
}

if ($lastStart -ge 0) {
    $matchRange = $d.Range($lastStart, $lastEnd)

    $targetPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $matchRange.Start -and $p.Range.End -ge $matchRange.End) {
            $targetPara = $p
        }
    }

    if ($targetPara -ne $null) {
        $targetPara.Range.Delete()
    }
}
